$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per cluster analysis rerun (2 variables, no SiO2)
$ws.Range("B2").Value = 56
$ws.Range("B3").Value = 43

# Remove the now-obsolete third data row (area index 2)
$ws.Range("A4:B4").Delete()
